# Restore C10 on the active sheet from 18 to 1 (value-only edit captured
# by the commit's "restore from revision" save).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C10").Value = 1
